$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9776
$ws1.Range("F4").Value = 36
$ws1.Range("F5").Value = 564
$ws1.Range("F6").Value = 471

# Sheet "全部类型" (All Types) - same underlying data, update matching rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9776
$ws4.Range("F4").Value = 36
$ws4.Range("F5").Value = 564
$ws4.Range("F7").Value = 471
